$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear (incl. the cell record itself) the A-column cells that no
# longer hold a date value in the new layout
$ws.Range("A3").Clear()
$ws.Range("A4").Clear()
$ws.Range("A5").Clear()
$ws.Range("A6").Clear()
$ws.Range("A8").Clear()
$ws.Range("A9").Clear()
$ws.Range("A10").Clear()
$ws.Range("A12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()
$ws.Range("A17").Clear()

# Set cells that keep/reuse existing shared strings (order doesn't affect final index assignment)
$ws.Range("A1").Value2 = "date"
$ws.Range("C1").Value2 = "approach"
$ws.Range("A2").Value2 = 1.16
$ws.Range("B2").Value2 = "project analysis"
$ws.Range("C2").Value2 = "work with the target site to develop a series of features needed to be implemented, came up with a list of features and general idea of approach"
$ws.Range("B3").Value2 = "visualization library select"
$ws.Range("C3").Value2 = "searched online and asked people, know there are popular libraries like d3, high charts and echarts"
$ws.Range("B4").Value2 = "try on d3"
$ws.Range("C4").Value2 = "went through a tutorial to see the basic logic of d3, found out that it is powerful, yet need to build everything from scratch, the learning curve is steep"
$ws.Range("B5").Value2 = "try on echarts"
$ws.Range("C5").Value2 = "went through  tutorial and demos to see the basic logic of echats and the effect it can achieve, found stack bar charts and line charts with gradient, also went through api, found out there are events and predefined functions like tooltip that can be used"
$ws.Range("B6").Value2 = "Decide to use echats"
$ws.Range("A7").Value2 = 1.17
$ws.Range("B7").Value2 = "decide data interface format"
$ws.Range("C7").Value2 = "went through tutorial and apis to see how echats deal with data, came up with data interface format: array of objs"
$ws.Range("B8").Value2 = "generate fake data for drawing bar chart"
$ws.Range("B9").Value2 = "followed examples to map data to graph"
$ws.Range("B10").Value2 = "explore possible approach for scroll bar"
$ws.Range("C10").Value2 = "after going through examples and apis, decided to use datazoom to implement scroll bar, set bar handle width to fixed, changed styling of the bar."
$ws.Range("D10").Value2 = "can't make the scrollbar identical to browzer's default scrollbar"
$ws.Range("A11").Value2 = 1.18
$ws.Range("B11").Value2 = "encountered prob that after setting the scrollbar, the stack bars on the graph starts to stack over eachother"
$ws.Range("C11").Value2 = "explored examples and apis, after trying all the parameters, found the filter parameter can be used to control this"
$ws.Range("B12").Value2 = "adjust chart layout"
$ws.Range("C12").Value2 = "explored axis apis and settings to adjust layout "
$ws.Range("D12").Value2 = "axis and text position can't be adjusted"
$ws.Range("B13").Value2 = "looking for ways to add click event for axis"
$ws.Range("C13").Value2 = "went through apis and examples, failed to make the axis response to click. Create buttons at the bottom instead to control  sorting"
$ws.Range("D13").Value2 = "adjust layout"
$ws.Range("B14").Value2 = "change the data input method"
$ws.Range("C14").Value2 = "found it not convinient to write data into series, learned dataset using reference, change data interface format into array of arrays"
$ws.Range("B15").Value2 = "ask whether the visulization is acceptable"
$ws.Range("A16").Value2 = 1.19
$ws.Range("B16").Value2 = "ajax dataloading"
$ws.Range("C16").Value2 = "decided to upload csv to github to use ajax request. "
$ws.Range("B17").Value2 = "dataparsing for chart3 and 4"
$ws.Range("A18").Value2 = 1.2
$ws.Range("B18").Value2 = "restructured Dom, added chart4 in"
$ws.Range("C19").Value2 = "found out that although I tried to make a copy of the origin data then process the copy, I was modifiying the original data, since the elements in the data array are objs. So I was passing reference to the copy."
$ws.Range("C20").Value2 = "since console.log() is asynchronous, it gets fired after the main body of the function is completed, the printed result is not the result at the expected moment"
$ws.Range("A21").Value2 = 1.21
$ws.Range("A25").Value2 = 1.22

# Set brand-new text cells in the same order the entries were introduced,
# so the rebuilt shared-strings table matches the target ordering
$ws.Range("B19").Value2 = "found bug in dataparsing for graph3 and 4"
$ws.Range("D19").Value2 = "to be fixed"
$ws.Range("C21").Value2 = "searched online, found a way of using JSON.parse(JSON.stringify(array)) to deep copy array of objs"
$ws.Range("B1").Value2 = "obj/prob encountered"
$ws.Range("B21").Value2 = "fixed data in graph3 and 4"
$ws.Range("B22").Value2 = "found duplication in rawdata, wonder whether need to deduplicate. "
$ws.Range("D22").Value2 = "dedup in the future"
$ws.Range("B23").Value2 = "fixed data fetching with fake data in chart2, all four charts are displayed with fake data now"
$ws.Range("B24").Value2 = "found bug with stacking feature in the library"
$ws.Range("C24").Value2 = "compared with demo code, tried different approaches to solve the prob. At last, got answer from the library team that it is the bug in the library that leads to the issue"
$ws.Range("D24").Value2 = "maybe need to change the way from fetching from dataset to extract data and fetch into series"
$ws.Range("D1").Value2 = "follow up"

# Row heights for the newly added rows (21-24); rows 1-20 keep their existing heights
$ws.Rows.Item(21).RowHeight = 43.2
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 57.6
$ws.Rows.Item(24).RowHeight = 72

# Update the sheet view selection to D4 (and drop the old scrolled
# topLeftCell/selection state)
$ws.Range("D4").Select()
